# Update "Latest HO Xliff Generate Date" / handoff/handback datetime stamps
# for the b512c95f-9fa1-403b-a53b-5c03f44ede5c row, as part of regenerating
# the handback report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for b512c95f row (row 4)
$overview.Range("G4").Value = "2016-09-02 22:50:46"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for b512c95f row (row 4)
$zhcn.Range("H4").Value = "2016-09-02 22:50:42"
$zhcn.Range("K4").Value = "2016-09-02 22:51:12"

# de-de sheet: Correspond Handoff Datetime (shared with Overview value) / Correspond Handback DateTime
$dede.Range("H4").Value = "2016-09-02 22:50:46"
$dede.Range("K4").Value = "2016-09-02 22:51:19"
